# Apply updated odds values to specific cells per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("M7").Value = 1.04
$ws.Range("N7").Value = 13
$ws.Range("O7").Value = 1.29
$ws.Range("P7").Value = 3.5
$ws.Range("Q7").Value = 1.88
$ws.Range("R7").Value = 1.93
$ws.Range("AJ7").Value = 26

# Row 8
$ws.Range("G8").Value = 2.55
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 3.25
$ws.Range("L8").Value = 3.75
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 7
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.73
$ws.Range("W8").Value = 7
$ws.Range("X8").Value = 11
$ws.Range("Y8").Value = 11
$ws.Range("Z8").Value = 26
$ws.Range("AA8").Value = 23
$ws.Range("AC8").Value = 7
$ws.Range("AD8").Value = 6
$ws.Range("AE8").Value = 17
$ws.Range("AI8").Value = 13
$ws.Range("AK8").Value = 29
$ws.Range("AU8").Value = 8.5
$ws.Range("AW8").Value = 4.75
$ws.Range("AX8").Value = 17
$ws.Range("AY8").Value = 29
$ws.Range("AZ8").Value = 51
$ws.Range("BB8").Value = 251

# Row 9
$ws.Range("G9").Value = 2.35
$ws.Range("H9").Value = 2.8
$ws.Range("J9").Value = 3.25
$ws.Range("M9").Value = 1.13
$ws.Range("N9").Value = 6
$ws.Range("O9").Value = 1.57
$ws.Range("P9").Value = 2.25
$ws.Range("Q9").Value = 2.88
$ws.Range("R9").Value = 1.4
$ws.Range("S9").Value = 1.62
$ws.Range("T9").Value = 2.2
$ws.Range("U9").Value = 2.25
$ws.Range("V9").Value = 1.57
$ws.Range("W9").Value = 5.5
$ws.Range("Y9").Value = 11
$ws.Range("Z9").Value = 23
$ws.Range("AC9").Value = 5.5
$ws.Range("AF9").Value = 81
$ws.Range("AH9").Value = 7.5
$ws.Range("AT9").Value = 2.2

# Row 10
$ws.Range("G10").Value = 1.45
$ws.Range("H10").Value = 4.1
$ws.Range("I10").Value = 7.5
$ws.Range("J10").Value = 2.05
$ws.Range("L10").Value = 7.5
$ws.Range("U10").Value = 2.38
$ws.Range("V10").Value = 1.53
$ws.Range("W10").Value = 5
$ws.Range("X10").Value = 6
$ws.Range("Y10").Value = 9
$ws.Range("Z10").Value = 9
$ws.Range("AD10").Value = 8
$ws.Range("AE10").Value = 23
$ws.Range("AH10").Value = 15
$ws.Range("AN10").Value = 3.2
$ws.Range("AO10").Value = 7.5
$ws.Range("AU10").Value = 10
$ws.Range("AW10").Value = 8.5

# Row 12
$ws.Range("I12").Value = 1.9
$ws.Range("J12").Value = 4.05
$ws.Range("L12").Value = 2.45
$ws.Range("S12").Value = 1.32
$ws.Range("T12").Value = 3.2
$ws.Range("V12").Value = 2.07
$ws.Range("W12").Value = 12.5
$ws.Range("X12").Value = 22
$ws.Range("Y12").Value = 12.5
$ws.Range("AA12").Value = 32
$ws.Range("AH12").Value = 8.25
$ws.Range("AI12").Value = 10
$ws.Range("AK12").Value = 17
$ws.Range("AL12").Value = 14
$ws.Range("AM12").Value = 22
$ws.Range("AO12").Value = 20
$ws.Range("AP12").Value = 24
$ws.Range("AT12").Value = 2.87
$ws.Range("AY12").Value = 16.5
$ws.Range("AZ12").Value = 32
$ws.Range("BA12").Value = 55
$ws.Range("BB12").Value = 175

# Row 19
$ws.Range("G19").Value = 1.6
$ws.Range("H19").Value = 3.6
$ws.Range("I19").Value = 5.4
$ws.Range("J19").Value = 2.1
$ws.Range("K19").Value = 2.18
$ws.Range("L19").Value = 5.5
$ws.Range("M19").Value = 1.03
$ws.Range("N19").Value = 12.6
$ws.Range("Q19").Value = 1.9
$ws.Range("R19").Value = 1.82
$ws.Range("W19").Value = 6.3
$ws.Range("X19").Value = 7.1
$ws.Range("Z19").Value = 11.5
$ws.Range("AA19").Value = 13
$ws.Range("AD19").Value = 7.1
$ws.Range("AE19").Value = 17
$ws.Range("AH19").Value = 14
$ws.Range("AI19").Value = 35
$ws.Range("AJ19").Value = 17.5
$ws.Range("AL19").Value = 60
$ws.Range("AM19").Value = 60
$ws.Range("AO19").Value = 7.5
$ws.Range("AP19").Value = 16.5
$ws.Range("AQ19").Value = 24
$ws.Range("AS19").Value = 200
$ws.Range("AT19").Value = 2.62
$ws.Range("AW19").Value = 7
$ws.Range("AX19").Value = 32
$ws.Range("AY19").Value = 35
$ws.Range("AZ19").Value = 200
$ws.Range("BB19").Value = 450
